$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 17, keeping only the header row and the first data row
$ws.Range("A3:B17").EntireRow.Delete()

# Update the remaining data row (row 2) with the new date and value
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 1.176843378132464
